$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2791963333333333
$ws.Range("H2").Value = 0.837589
$ws.Range("I2").Value = 0.008912184157424861
$ws.Range("J2").Value = 0.009090659364840875
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.343285
$ws.Range("N2").Value = 1.029855
$ws.Range("O2").Value = 0.9867450936054706
$ws.Range("P2").Value = 0.9867450936054707
$ws.Range("Q2").Value = 0.09584391328833333
$ws.Range("R2").Value = 0.862595219595
$ws.Range("S2").Value = 0.008794053990647387
$ws.Range("T2").Value = 0.008970163525895358
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2791963333333333
$ws.Range("H3").Value = 0.837589
$ws.Range("I3").Value = 0.008912184157424861
$ws.Range("J3").Value = 0.009090659364840875
$ws.Range("M3").Value = 0.004611333333333334
$ws.Range("N3").Value = 0.013834
$ws.Range("O3").Value = 0.0132549063945294
$ws.Range("P3").Value = 0.01325490639452941
$ws.Range("Q3").Value = 0.001287467358444445
$ws.Range("R3").Value = 0.011587206226
$ws.Range("S3").Value = 0.0001181301667774744
$ws.Range("T3").Value = 0.000120495838945518
$ws.Range("G4").Value = 2.510701
$ws.Range("H4").Value = 7.532103
$ws.Range("I4").Value = 0.08014370894160773
$ws.Range("J4").Value = 0.0817486651256118
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.343285
$ws.Range("N4").Value = 1.029855
$ws.Range("O4").Value = 0.9867450936054706
$ws.Range("P4").Value = 0.9867450936054707
$ws.Range("Q4").Value = 0.861885992785
$ws.Range("R4").Value = 7.756973935065
$ws.Range("S4").Value = 0.07908141158147632
$ws.Range("T4").Value = 0.08066509422149409
$ws.Range("G5").Value = 2.510701
$ws.Range("H5").Value = 7.532103
$ws.Range("I5").Value = 0.08014370894160773
$ws.Range("J5").Value = 0.0817486651256118
$ws.Range("M5").Value = 0.004611333333333334
$ws.Range("N5").Value = 0.013834
$ws.Range("O5").Value = 0.0132549063945294
$ws.Range("P5").Value = 0.01325490639452941
$ws.Range("Q5").Value = 0.01157767921133333
$ws.Range("R5").Value = 0.104199112902
$ws.Range("S5").Value = 0.00106229736013142
$ws.Range("T5").Value = 0.001083570904117715
$ws.Range("G6").Value = 16.644438
$ws.Range("H6").Value = 49.933314
$ws.Range("I6").Value = 0.5313046016107197
$ws.Range("J6").Value = 0.5419444960853593
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.343285
$ws.Range("N6").Value = 1.029855
$ws.Range("O6").Value = 0.9867450936054706
$ws.Range("P6").Value = 0.9867450936054707
$ws.Range("Q6").Value = 5.71378589883
$ws.Range("R6").Value = 51.42407308947
$ws.Range("S6").Value = 0.5242622088493869
$ws.Range("T6").Value = 0.5347610725187175
$ws.Range("G7").Value = 16.644438
$ws.Range("H7").Value = 49.933314
$ws.Range("I7").Value = 0.5313046016107197
$ws.Range("J7").Value = 0.5419444960853593
$ws.Range("M7").Value = 0.004611333333333334
$ws.Range("N7").Value = 0.013834
$ws.Range("O7").Value = 0.0132549063945294
$ws.Range("P7").Value = 0.01325490639452941
$ws.Range("Q7").Value = 0.076753051764
$ws.Range("R7").Value = 0.6907774658760001
$ws.Range("S7").Value = 0.007042392761332826
$ws.Range("T7").Value = 0.007183423566641847
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.84514
$ws.Range("H8").Value = 3.69028
$ws.Range("I8").Value = 0.0588984363795283
$ws.Range("J8").Value = 0.04005195679609568
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.343285
$ws.Range("N8").Value = 1.029855
$ws.Range("O8").Value = 0.9867450936054706
$ws.Range("P8").Value = 0.9867450936054707
$ws.Range("Q8").Value = 0.6334088849
$ws.Range("R8").Value = 3.8004533094
$ws.Range("S8").Value = 0.05811774311853351
$ws.Range("T8").Value = 0.0395210718578457
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.84514
$ws.Range("H9").Value = 3.69028
$ws.Range("I9").Value = 0.0588984363795283
$ws.Range("J9").Value = 0.04005195679609568
$ws.Range("M9").Value = 0.004611333333333334
$ws.Range("N9").Value = 0.013834
$ws.Range("O9").Value = 0.0132549063945294
$ws.Range("P9").Value = 0.01325490639452941
$ws.Range("Q9").Value = 0.008508555586666668
$ws.Range("R9").Value = 0.05105133352000001
$ws.Range("S9").Value = 0.000780693260994793
$ws.Range("T9").Value = 0.0005308849382499842
$ws.Range("G10").Value = 10.04801166666667
$ws.Range("H10").Value = 30.144035
$ws.Range("I10").Value = 0.3207410689107194
$ws.Range("J10").Value = 0.3271642226280922
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.343285
$ws.Range("N10").Value = 1.029855
$ws.Range("O10").Value = 0.9867450936054706
$ws.Range("P10").Value = 0.9867450936054707
$ws.Range("Q10").Value = 3.449331684991667
$ws.Range("R10").Value = 31.043985164925
$ws.Range("S10").Value = 0.3164896760654265
$ws.Range("T10").Value = 0.322827691481518
$ws.Range("G11").Value = 10.04801166666667
$ws.Range("H11").Value = 30.144035
$ws.Range("I11").Value = 0.3207410689107194
$ws.Range("J11").Value = 0.3271642226280922
$ws.Range("M11").Value = 0.004611333333333334
$ws.Range("N11").Value = 0.013834
$ws.Range("O11").Value = 0.0132549063945294
$ws.Range("P11").Value = 0.01325490639452941
$ws.Range("Q11").Value = 0.04633473113222223
$ws.Range("R11").Value = 0.41701258019
$ws.Range("S11").Value = 0.004251392845292891
$ws.Range("T11").Value = 0.004336531146574342
